# Daily scrape refresh: the sheet now reflects a single (updated) opportunity
# record in row 2, and the previously-scraped rows 3-7 are gone. A couple of
# column widths were also tweaked to better fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the now-stale rows 3-7 (only the header + one data row remain) ---
$ws.Range("A3:H7").EntireRow.Delete()

# --- Re-fit a handful of column widths (raw OOXML width = ColumnWidth + 5/6) ---
$ws.Columns.Item(3).ColumnWidth = 30.166666666666668   # C: 102 -> 31
$ws.Columns.Item(4).ColumnWidth = 100.16666666666667   # D: 57  -> 101
$ws.Columns.Item(6).ColumnWidth = 14.166666666666666   # F: 16  -> 15
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666   # G: 16  -> 15
$ws.Columns.Item(8).ColumnWidth = 23.166666666666668   # H: 57  -> 24

# --- Refresh the single remaining data row (row 2) with the new scrape ---
$a2 = $ws.Range("A2")
$a2.NumberFormat = "@"          # force the numeric-looking ID to stay text
$a2.Value = "1327883"
$a2.Style = "Normal"            # drop the temporary text format again

$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327883"
$ws.Range("C2").Value = "Mobile Application Developer"
$ws.Range("D2").Value = "El-Mahalla El-Kubra, Al Mahalah Al Kubra (Part 2), El Mahalla El Kubra, Gharbia Governorate, Egypt"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "8 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Positive Kids academy"
